$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D13").Value = "2016-03-09 16:21:47"
$wsZh.Range("D14").Value = "2016-03-09 16:21:47"

$wsDe.Range("D13").Value = "2016-03-09 16:21:53"
$wsDe.Range("D14").Value = "2016-03-09 16:21:53"
